$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 50 <-> row 51 (columns F:V); A:E unchanged ---
$ws.Range("F50").Value = 'TSC'
$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 'Zeleznicar Pancevo'
$ws.Range("I50").Value = 3
$ws.Range("J50").Value = 1.21
$ws.Range("K50").Value = '14/09/2023 09:13'
$ws.Range("L50").Value = 1.33
$ws.Range("M50").Value = '16/09/2023 18:52'
$ws.Range("N50").Value = 5.61
$ws.Range("O50").Value = '14/09/2023 09:13'
$ws.Range("P50").Value = 4.42
$ws.Range("Q50").Value = '16/09/2023 18:52'
$ws.Range("R50").Value = 9.23
$ws.Range("S50").Value = '14/09/2023 09:13'
$ws.Range("T50").Value = 10.82
$ws.Range("U50").Value = '16/09/2023 18:52'
$ws.Range("V50").Value = 'https://www.betexplorer.com/football/serbia/super-liga/tsc-backa-topola-zeleznicar-pancevo/xOIdSqWO/'

$ws.Range("F51").Value = 'Vozdovac'
$ws.Range("G51").Value = 3
$ws.Range("H51").Value = 'Mladost'
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = 2.07
$ws.Range("K51").Value = '14/09/2023 09:13'
$ws.Range("L51").Value = 1.83
$ws.Range("M51").Value = '16/09/2023 18:52'
$ws.Range("N51").Value = 3.13
$ws.Range("O51").Value = '14/09/2023 09:13'
$ws.Range("P51").Value = 3.58
$ws.Range("Q51").Value = '16/09/2023 18:52'
$ws.Range("R51").Value = 3.2
$ws.Range("S51").Value = '14/09/2023 09:13'
$ws.Range("T51").Value = 4.06
$ws.Range("U51").Value = '16/09/2023 18:52'
$ws.Range("V51").Value = 'https://www.betexplorer.com/football/serbia/super-liga/fk-vozdovac-mladost-lucani/25QqVon6/'

# --- Swap row 90 <-> row 91 (columns F:V); A:E unchanged ---
$ws.Range("F90").Value = 'Radnicki Nis'
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 'Sp. Subotica'
$ws.Range("I90").Value = 1
$ws.Range("J90").Value = 1.7
$ws.Range("K90").Value = '27/10/2023 06:42'
$ws.Range("L90").Value = 1.66
$ws.Range("M90").Value = '28/10/2023 18:23'
$ws.Range("N90").Value = 3.42
$ws.Range("O90").Value = '27/10/2023 06:42'
$ws.Range("P90").Value = 3.65
$ws.Range("Q90").Value = '28/10/2023 18:23'
$ws.Range("R90").Value = 4.23
$ws.Range("S90").Value = '27/10/2023 06:42'
$ws.Range("T90").Value = 5.08
$ws.Range("U90").Value = '28/10/2023 18:23'
$ws.Range("V90").Value = 'https://www.betexplorer.com/football/serbia/super-liga/radnicki-nis-spartak-subotica/2qDshl5f/'

$ws.Range("F91").Value = 'IMT Novi Beograd'
$ws.Range("G91").Value = 1
$ws.Range("H91").Value = 'Crvena zvezda'
$ws.Range("I91").Value = 2
$ws.Range("J91").Value = 8.15
$ws.Range("K91").Value = '27/10/2023 06:42'
$ws.Range("L91").Value = 24.2
$ws.Range("M91").Value = '28/10/2023 18:29'
$ws.Range("N91").Value = 5.6
$ws.Range("O91").Value = '27/10/2023 06:42'
$ws.Range("P91").Value = 9.529999999999999
$ws.Range("Q91").Value = '28/10/2023 18:29'
$ws.Range("R91").Value = 1.23
$ws.Range("S91").Value = '27/10/2023 06:42'
$ws.Range("T91").Value = 1.09
$ws.Range("U91").Value = '28/10/2023 18:21'
$ws.Range("V91").Value = 'https://www.betexplorer.com/football/serbia/super-liga/imt-novi-beograd-crvena-zvezda/SjAgknkD/'

# --- Append new rows 100, 101, 102 (copy A/E formatting from row 99) ---
# Row 100
$ws.Range("A99").Copy()
$ws.Range("A100").PasteSpecial(-4122)
$ws.Range("E99").Copy()
$ws.Range("E100").PasteSpecial(-4122)
$ws.Range("A100").Value = 99
$ws.Range("B100").Value = 'serbia'
$ws.Range("C100").Value = 'super-liga'
$ws.Range("D100").Value = '2023-2024'
$ws.Range("E100").Value = 45236.70833333334
$ws.Range("F100").Value = 'Sp. Subotica'
$ws.Range("G100").Value = 1
$ws.Range("H100").Value = 'Radnicki 1923'
$ws.Range("I100").Value = 3
$ws.Range("J100").Value = 2.13
$ws.Range("K100").Value = '02/11/2023 11:12'
$ws.Range("L100").Value = 2.67
$ws.Range("M100").Value = '06/11/2023 16:59'
$ws.Range("N100").Value = 3.1
$ws.Range("O100").Value = '02/11/2023 11:12'
$ws.Range("P100").Value = 3.27
$ws.Range("Q100").Value = '06/11/2023 16:59'
$ws.Range("R100").Value = 3.1
$ws.Range("S100").Value = '02/11/2023 11:12'
$ws.Range("T100").Value = 2.56
$ws.Range("U100").Value = '06/11/2023 16:59'
$ws.Range("V100").Value = 'https://www.betexplorer.com/football/serbia/super-liga/spartak-subotica-radnicki-1923/zcXK57Zg/'

# Row 101
$ws.Range("A99").Copy()
$ws.Range("A101").PasteSpecial(-4122)
$ws.Range("E99").Copy()
$ws.Range("E101").PasteSpecial(-4122)
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = 'serbia'
$ws.Range("C101").Value = 'super-liga'
$ws.Range("D101").Value = '2023-2024'
$ws.Range("E101").Value = 45236.77083333334
$ws.Range("F101").Value = 'Napredak'
$ws.Range("G101").Value = 0
$ws.Range("H101").Value = 'IMT Novi Beograd'
$ws.Range("I101").Value = 4
$ws.Range("J101").Value = 2.24
$ws.Range("K101").Value = '02/11/2023 08:13'
$ws.Range("L101").Value = 2.3
$ws.Range("M101").Value = '06/11/2023 18:26'
$ws.Range("N101").Value = 3.03
$ws.Range("O101").Value = '02/11/2023 08:13'
$ws.Range("P101").Value = 3.37
$ws.Range("Q101").Value = '06/11/2023 18:26'
$ws.Range("R101").Value = 2.95
$ws.Range("S101").Value = '02/11/2023 08:13'
$ws.Range("T101").Value = 2.94
$ws.Range("U101").Value = '06/11/2023 18:26'
$ws.Range("V101").Value = 'https://www.betexplorer.com/football/serbia/super-liga/napredak-imt-novi-beograd/vVwC7T4s/'

# Row 102
$ws.Range("A99").Copy()
$ws.Range("A102").PasteSpecial(-4122)
$ws.Range("E99").Copy()
$ws.Range("E102").PasteSpecial(-4122)
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = 'serbia'
$ws.Range("C102").Value = 'super-liga'
$ws.Range("D102").Value = '2023-2024'
$ws.Range("E102").Value = 45236.79166666666
$ws.Range("F102").Value = 'Mladost'
$ws.Range("G102").Value = 1
$ws.Range("H102").Value = 'Javor'
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 2.26
$ws.Range("K102").Value = '02/11/2023 08:13'
$ws.Range("L102").Value = 2.28
$ws.Range("M102").Value = '06/11/2023 18:56'
$ws.Range("N102").Value = 3.04
$ws.Range("O102").Value = '02/11/2023 08:13'
$ws.Range("P102").Value = 3.17
$ws.Range("Q102").Value = '06/11/2023 18:56'
$ws.Range("R102").Value = 2.91
$ws.Range("S102").Value = '02/11/2023 08:13'
$ws.Range("T102").Value = 3.15
$ws.Range("U102").Value = '06/11/2023 18:56'
$ws.Range("V102").Value = 'https://www.betexplorer.com/football/serbia/super-liga/mladost-lucani-javor/d6YG6mKm/'

$excel.CutCopyMode = $false
